# Update the date line and all the division problems in the table.
# "531÷5=" appears twice in the original document with two different
# replacement values, so those two cells are addressed directly via the
# Tables collection (row 5 col 4, and row 17 col 5, 1-based) instead of a
# global Find/Replace.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Date header
Replace-Text "2024-08-17 Saturday" "2024-08-18 Sunday"

# Row 1 (table row 1)
Replace-Text "676÷7=" "103÷4="
Replace-Text "863÷4=" "845÷7="
Replace-Text "563÷9=" "744÷3="
Replace-Text "146÷4=" "394÷8="
Replace-Text "756÷8=" "728÷3="

# Row 2 (table row 5) - note: "531÷5=" is ambiguous, handled below via Cell access
Replace-Text "520÷4=" "515÷7="
Replace-Text "868÷3=" "750÷6="
Replace-Text "551÷5=" "906÷6="
Replace-Text "839÷8=" "583÷2="

# Row 3 (table row 9)
Replace-Text "487÷8=" "897÷3="
Replace-Text "587÷6=" "958÷8="
Replace-Text "584÷8=" "167÷9="
Replace-Text "943÷5=" "112÷4="
Replace-Text "980÷8=" "128÷3="

# Row 4 (table row 13)
Replace-Text "664÷4=" "882÷2="
Replace-Text "879÷8=" "680÷4="
Replace-Text "718÷3=" "281÷2="
Replace-Text "773÷5=" "600÷4="
Replace-Text "669÷7=" "653÷5="

# Row 5 (table row 17) - note: "531÷5=" is ambiguous, handled below via Cell access
Replace-Text "272÷7=" "609÷2="
Replace-Text "684÷4=" "245÷2="
Replace-Text "312÷5=" "711÷6="
Replace-Text "168÷4=" "722÷6="

# The two remaining "531÷5=" cells (table row 5 col 4, and table row 17 col 5)
# get different replacements, so address the exact table cells directly.
$table = $d.Tables.Item(1)
$table.Cell(5, 4).Range.Text = "879÷2="
$table.Cell(17, 5).Range.Text = "783÷8="
